# Add a new book entry ("prueba" / "otr") to the wishlist as row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "prueba"
$ws.Range("B2").Value = "otr"

# Touch C2 (no Editorial value for this entry) so it is materialized as an
# empty cell in the sheet, keeping it in line with A2:B2 and extending the
# used range to A1:C2, without picking up a new/explicit cell style.
$ws.Range("C2").Style = "Normal"
